# Auto-generated edit script: updates pricing/profit columns (H-N)
# across multiple worksheets to reflect refreshed market data.
$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 630.2
$ws.Cells.Item(9, 9).Value = 287.75
$ws.Cells.Item(9, 11).Value = 287.75
$ws.Cells.Item(9, 13).Value = -118.75
$ws.Cells.Item(17, 8).Value = 2993.5
$ws.Cells.Item(17, 10).Value = 2993.5
$ws.Cells.Item(17, 12).Value = 8980.5
$ws.Cells.Item(17, 14).Value = -9316.5
$ws.Cells.Item(92, 8).Value = 1029.3334
$ws.Cells.Item(92, 9).Value = 1044
$ws.Cells.Item(92, 10).Value = 1000
$ws.Cells.Item(92, 11).Value = 1044
$ws.Cells.Item(92, 12).Value = 1000
$ws.Cells.Item(92, 13).Value = 204
$ws.Cells.Item(92, 14).Value = -3496
$ws.Cells.Item(116, 8).Value = 5915.8
$ws.Cells.Item(116, 9).Value = 6169.75
$ws.Cells.Item(116, 11).Value = 6169.75
$ws.Cells.Item(116, 13).Value = -2727.75
$ws.Cells.Item(138, 8).Value = 3776.5898
$ws.Cells.Item(138, 10).Value = 3278.862
$ws.Cells.Item(138, 12).Value = 9836.585999999999
$ws.Cells.Item(138, 14).Value = -20116.586

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6162.143
$ws.Cells.Item(32, 9).Value = 4036.2258
$ws.Cells.Item(32, 11).Value = 4036.2258
$ws.Cells.Item(32, 13).Value = -3749.2258
$ws.Cells.Item(61, 8).Value = 4998.5
$ws.Cells.Item(61, 10).Value = 4998
$ws.Cells.Item(61, 12).Value = 4998
$ws.Cells.Item(61, 14).Value = -5422
$ws.Cells.Item(97, 8).Value = 813.381
$ws.Cells.Item(97, 10).Value = 1283
$ws.Cells.Item(97, 12).Value = 1283
$ws.Cells.Item(97, 14).Value = -2275
$ws.Cells.Item(102, 8).Value = 4204.636
$ws.Cells.Item(102, 9).Value = 4535.1
$ws.Cells.Item(102, 10).Value = 900
$ws.Cells.Item(102, 11).Value = 4535.1
$ws.Cells.Item(102, 12).Value = 900
$ws.Cells.Item(102, 13).Value = -2913.1
$ws.Cells.Item(102, 14).Value = -4144
$ws.Cells.Item(132, 8).Value = 2898.0488
$ws.Cells.Item(132, 9).Value = 2898.0488
$ws.Cells.Item(132, 11).Value = 8694.1464
$ws.Cells.Item(132, 13).Value = -6164.1464
$ws.Cells.Item(136, 8).Value = 4998.5
$ws.Cells.Item(136, 10).Value = 4998
$ws.Cells.Item(136, 12).Value = 14994
$ws.Cells.Item(136, 14).Value = -20094

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(16, 8).Value = 675
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 14).ClearContents()
$ws.Cells.Item(20, 8).Value = 922.2857
$ws.Cells.Item(20, 9).Value = 913.9167
$ws.Cells.Item(20, 10).Value = 972.5
$ws.Cells.Item(20, 11).Value = 913.9167
$ws.Cells.Item(20, 12).Value = 972.5
$ws.Cells.Item(20, 13).Value = -666.9167
$ws.Cells.Item(20, 14).Value = -1466.5
$ws.Cells.Item(86, 8).Value = 15173585
$ws.Cells.Item(86, 9).Value = 25334.416
$ws.Cells.Item(86, 10).Value = 33351486
$ws.Cells.Item(86, 11).Value = 25334.416
$ws.Cells.Item(86, 12).Value = 33351486
$ws.Cells.Item(86, 13).Value = -24211.416
$ws.Cells.Item(86, 14).Value = -33353732
$ws.Cells.Item(89, 8).Value = 15173585
$ws.Cells.Item(89, 9).Value = 25334.416
$ws.Cells.Item(89, 10).Value = 33351486
$ws.Cells.Item(89, 11).Value = 126672.08
$ws.Cells.Item(89, 12).Value = 166757430
$ws.Cells.Item(89, 13).Value = -121056.08
$ws.Cells.Item(89, 14).Value = -166768662
$ws.Cells.Item(94, 8).Value = 1399.4
$ws.Cells.Item(94, 9).Value = 1399.4
$ws.Cells.Item(94, 11).Value = 1399.4
$ws.Cells.Item(94, 13).Value = -948.4000000000001

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 50000180
$ws.Cells.Item(7, 9).Value = 55555744
$ws.Cells.Item(7, 11).Value = 55555744
$ws.Cells.Item(7, 13).Value = -55555631
$ws.Cells.Item(31, 8).Value = 2413.9583
$ws.Cells.Item(31, 9).Value = 1896.1666
$ws.Cells.Item(31, 11).Value = 1896.1666
$ws.Cells.Item(31, 13).Value = -1601.1666
$ws.Cells.Item(34, 8).Value = 2413.9583
$ws.Cells.Item(34, 9).Value = 1896.1666
$ws.Cells.Item(34, 11).Value = 1896.1666
$ws.Cells.Item(34, 13).Value = -1694.1666
$ws.Cells.Item(99, 8).Value = 12893.111
$ws.Cells.Item(99, 9).Value = 5819.25
$ws.Cells.Item(99, 11).Value = 5819.25
$ws.Cells.Item(99, 13).Value = -4321.25
$ws.Cells.Item(107, 8).Value = 244.3
$ws.Cells.Item(107, 9).Value = 223.28572
$ws.Cells.Item(107, 10).Value = 293.33334
$ws.Cells.Item(107, 11).Value = 223.28572
$ws.Cells.Item(107, 12).Value = 293.33334
$ws.Cells.Item(107, 13).Value = 1696.71428
$ws.Cells.Item(107, 14).Value = -4133.33334
$ws.Cells.Item(126, 8).Value = 12893.111
$ws.Cells.Item(126, 9).Value = 5819.25
$ws.Cells.Item(126, 11).Value = 17457.75
$ws.Cells.Item(126, 13).Value = -14987.75

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 571464.8
$ws.Cells.Item(2, 9).Value = 59
$ws.Cells.Item(2, 10).Value = 769259.1
$ws.Cells.Item(2, 11).Value = 354
$ws.Cells.Item(2, 12).Value = 4615554.6
$ws.Cells.Item(2, 13).Value = -241
$ws.Cells.Item(2, 14).Value = -4615780.6
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 14).ClearContents()
$ws.Cells.Item(23, 8).Value = 236.82353
$ws.Cells.Item(23, 9).Value = 129
$ws.Cells.Item(23, 10).Value = 281.75
$ws.Cells.Item(23, 11).Value = 387
$ws.Cells.Item(23, 12).Value = 845.25
$ws.Cells.Item(23, 13).Value = -152
$ws.Cells.Item(23, 14).Value = -1315.25
$ws.Cells.Item(27, 8).Value = 0
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 12).Value = 0
$ws.Cells.Item(27, 14).ClearContents()
$ws.Cells.Item(75, 8).Value = 1685.8334
$ws.Cells.Item(75, 10).Value = 1000
$ws.Cells.Item(75, 12).Value = 3000
$ws.Cells.Item(75, 14).Value = -4996
$ws.Cells.Item(78, 8).Value = 1685.8334
$ws.Cells.Item(78, 10).Value = 1000
$ws.Cells.Item(78, 12).Value = 9000
$ws.Cells.Item(78, 14).Value = -18984
$ws.Cells.Item(131, 8).Value = 1964.0625
$ws.Cells.Item(131, 10).Value = 1761.6666
$ws.Cells.Item(131, 12).Value = 5284.9998
$ws.Cells.Item(131, 14).Value = -15364.9998
$ws.Cells.Item(134, 8).Value = 15500
$ws.Cells.Item(134, 9).Value = 15500
$ws.Cells.Item(134, 11).Value = 46500
$ws.Cells.Item(134, 13).Value = -41430

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 185.88889
$ws.Cells.Item(2, 9).Value = 207.05556
$ws.Cells.Item(2, 10).Value = 143.55556
$ws.Cells.Item(2, 11).Value = 207.05556
$ws.Cells.Item(2, 12).Value = 143.55556
$ws.Cells.Item(2, 13).Value = -94.05556000000001
$ws.Cells.Item(2, 14).Value = -369.55556
$ws.Cells.Item(9, 8).Value = 3799.5
$ws.Cells.Item(9, 9).Value = 1000
$ws.Cells.Item(9, 10).Value = 6599
$ws.Cells.Item(9, 11).Value = 1000
$ws.Cells.Item(9, 12).Value = 6599
$ws.Cells.Item(9, 13).Value = -830
$ws.Cells.Item(9, 14).Value = -6939
$ws.Cells.Item(97, 8).Value = 1431.2941
$ws.Cells.Item(97, 9).Value = 1431.2941
$ws.Cells.Item(97, 11).Value = 1431.2941
$ws.Cells.Item(97, 13).Value = -935.2941000000001
$ws.Cells.Item(112, 8).Value = 92500
$ws.Cells.Item(112, 10).Value = 92500
$ws.Cells.Item(112, 12).Value = 92500
$ws.Cells.Item(112, 14).Value = -94716
$ws.Cells.Item(132, 8).Value = 23811814
$ws.Cells.Item(132, 9).Value = 2507.4167
$ws.Cells.Item(132, 10).Value = 166667650
$ws.Cells.Item(132, 11).Value = 7522.250100000001
$ws.Cells.Item(132, 12).Value = 500002950
$ws.Cells.Item(132, 13).Value = -4992.250100000001
$ws.Cells.Item(132, 14).Value = -500008010
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 10).Value = 0
$ws.Cells.Item(133, 12).Value = 0
$ws.Cells.Item(133, 14).ClearContents()

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(9, 8).Value = 663.6
$ws.Cells.Item(9, 9).Value = 329.5
$ws.Cells.Item(9, 10).Value = 2000
$ws.Cells.Item(9, 11).Value = 329.5
$ws.Cells.Item(9, 12).Value = 2000
$ws.Cells.Item(9, 13).Value = -105.5
$ws.Cells.Item(9, 14).Value = -2448
$ws.Cells.Item(16, 8).Value = 497.63635
$ws.Cells.Item(16, 9).Value = 497.63635
$ws.Cells.Item(16, 11).Value = 497.63635
$ws.Cells.Item(16, 13).Value = -327.63635
$ws.Cells.Item(61, 8).Value = 2896.125
$ws.Cells.Item(61, 9).Value = 3076.8462
$ws.Cells.Item(61, 10).Value = 2113
$ws.Cells.Item(61, 11).Value = 3076.8462
$ws.Cells.Item(61, 12).Value = 2113
$ws.Cells.Item(61, 13).Value = -2874.8462
$ws.Cells.Item(61, 14).Value = -2517
$ws.Cells.Item(82, 8).Value = 1139.6428
$ws.Cells.Item(82, 9).Value = 789.1667
$ws.Cells.Item(82, 10).Value = 1402.5
$ws.Cells.Item(82, 11).Value = 789.1667
$ws.Cells.Item(82, 12).Value = 1402.5
$ws.Cells.Item(82, 13).Value = -428.1667
$ws.Cells.Item(82, 14).Value = -2124.5
$ws.Cells.Item(85, 8).Value = 1139.6428
$ws.Cells.Item(85, 9).Value = 789.1667
$ws.Cells.Item(85, 10).Value = 1402.5
$ws.Cells.Item(85, 11).Value = 789.1667
$ws.Cells.Item(85, 12).Value = 1402.5
$ws.Cells.Item(85, 13).Value = 458.8333
$ws.Cells.Item(85, 14).Value = -3898.5
$ws.Cells.Item(113, 8).Value = 2896.125
$ws.Cells.Item(113, 9).Value = 3076.8462
$ws.Cells.Item(113, 10).Value = 2113
$ws.Cells.Item(113, 11).Value = 3076.8462
$ws.Cells.Item(113, 12).Value = 2113
$ws.Cells.Item(113, 13).Value = -906.8462
$ws.Cells.Item(113, 14).Value = -6453
$ws.Cells.Item(122, 8).Value = 1970
$ws.Cells.Item(122, 9).Value = 1970
$ws.Cells.Item(122, 11).Value = 5910
$ws.Cells.Item(122, 13).Value = -3460
$ws.Cells.Item(136, 8).Value = 4220.778
$ws.Cells.Item(136, 9).Value = 4220.778
$ws.Cells.Item(136, 11).Value = 12662.334
$ws.Cells.Item(136, 13).Value = -10112.334

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 2557
$ws.Cells.Item(96, 9).Value = 2730.111
$ws.Cells.Item(96, 10).Value = 999
$ws.Cells.Item(96, 11).Value = 2730.111
$ws.Cells.Item(96, 12).Value = 999
$ws.Cells.Item(96, 13).Value = -1357.111
$ws.Cells.Item(96, 14).Value = -3745
$ws.Cells.Item(132, 8).Value = 142858720
$ws.Cells.Item(132, 9).Value = 1804.2
$ws.Cells.Item(132, 11).Value = 5412.6
$ws.Cells.Item(132, 13).Value = -2882.6
$ws.Cells.Item(136, 8).Value = 8502.223
$ws.Cells.Item(136, 9).Value = 8838.823
$ws.Cells.Item(136, 11).Value = 26516.469
$ws.Cells.Item(136, 13).Value = -23966.469

